$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Locate the two paragraphs this edit revolves around, by content rather
# than by hard-coded index, so the script is resilient to the document's
# exact paragraph numbering.
# --------------------------------------------------------------------------
$metaParaIndex = -1
$imagePromptIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($metaParaIndex -eq -1 -and $t.StartsWith("Meta description")) {
        $metaParaIndex = $i
    }
    if ($t.StartsWith("Create a feature image for Electric Avenue")) {
        $imagePromptIndex = $i
    }
}

# --------------------------------------------------------------------------
# 1. Capture the "Meta description" paragraph's OOXML (run structure) before
#    removing it, so its exact formatting (leading empty run + bold run)
#    can be reused later for the new "Play Electric Avenue..." paragraph
#    that gets added near the end of the document.
# --------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item($metaParaIndex)
$metaXml = $metaPara.Range.WordOpenXML
# Strip the paragraph/revision-tracking identifiers Word stamps onto the
# captured OOXML so the re-inserted paragraph doesn't pick up bogus ids
# that were never present in the original document.
$cleanXml = $metaXml -replace ' w14:paraId="[^"]*"', '' `
                      -replace ' w14:textId="[^"]*"', '' `
                      -replace ' w:rsidR="[^"]*"', '' `
                      -replace ' w:rsidRDefault="[^"]*"', ''

# --------------------------------------------------------------------------
# 2. Insert a copy of that paragraph immediately before the final
#    "Create a feature image..." paragraph.
# --------------------------------------------------------------------------
$imagePara = $d.Paragraphs.Item($imagePromptIndex)
$insertPoint = $d.Range($imagePara.Range.Start, $imagePara.Range.Start)
$insertPoint.InsertXML($cleanXml)

# InsertXML brings along a trailing empty paragraph mark (an artifact of the
# captured range's own paragraph mark) - remove it. It now sits right before
# the (shifted) image-prompt paragraph.
$imagePara = $d.Paragraphs.Item($imagePromptIndex + 1)
$spuriousPara = $imagePara.Previous
$spuriousRange = $spuriousPara.Range
if (($spuriousRange.End - $spuriousRange.Start) -le 1) {
    $spuriousRange.Delete()
}

# --------------------------------------------------------------------------
# 3. Remove the original "Meta description" paragraph entirely (including
#    its paragraph mark).
# --------------------------------------------------------------------------
$d.Paragraphs.Item($metaParaIndex).Range.Delete()

# --------------------------------------------------------------------------
# 4. Rewrite the text of the paragraph pasted in step 2 (which currently
#    reads "Meta description: Read our review...") down to just the new
#    bold title text.
# --------------------------------------------------------------------------
$newTitleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("Meta description")) {
        $newTitleIndex = $i
        break
    }
}
$newTitlePara = $d.Paragraphs.Item($newTitleIndex)
$titleRange = $newTitlePara.Range
$titleTextOnly = $d.Range($titleRange.Start, $titleRange.End - 1)
$titleTextOnly.Text = "Play Electric Avenue Slot Free | Microgaming Review"

# --------------------------------------------------------------------------
# 5. Replace the final paragraph's (italic) text with the text that used to
#    be the meta description.
# --------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a feature image for Electric Avenue that features a happy Maya warrior with glasses in a cartoon style. The Maya warrior should be holding a boombox and standing in front of a neon background, with symbols from the game such as sneakers and Rubik's Cubes surrounding them. The image should convey the fun and energetic atmosphere of the game and incorporate the 80s theme, while also highlighting the potential for big wins through the use of multipliers and Wild Reels. Make sure the image is eye-catching and appeals to players who enjoy high variance slots with unique themes and features.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Electric Avenue slot from Microgaming and play for free at top online casinos. Features, pros, cons, and RTP information included.",
    2
)
